$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "CASE WHEN..." column (old column B); this shifts C:G left to B:F
$ws.Columns.Item(2).Delete()

# Update identifiers in row 3 (measure/dimension qualifiers)
$ws.Range("B3").Value = "iaest-dimension:estrato"
$ws.Range("D3").Value = "sdmx-dimension:refArea"
$ws.Range("E3").Value = "iaest-dimension:mes-y-ano"

# Update row 4 ("medida"/"dim" marker row)
$ws.Range("B4").Value = "dim"
$ws.Range("D4").Value = "dim"
$ws.Range("E4").Value = "dim"

# Update row 5 (data type row)
$ws.Range("B5").Value = "skos:Concept"
$ws.Range("D5").Value = "URI-Provincia"

# New row 6, re-using formatting from an existing styled cell
$ws.Range("A5").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("B6").Value = "mapping-estrato.xlsx"
